# Agregando un nuevo permiso
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Find first empty row after the existing data (row 58 in this case)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = $newRow - 1
$ws.Cells.Item($newRow, 2).Value = "dashboard.list"
$ws.Cells.Item($newRow, 3).Value = "Permite visualizar la pantalla de dashboard"

# Match the style used by the rest of column B for this new row
$ws.Cells.Item($newRow - 1, 2).Copy()
$ws.Cells.Item($newRow, 2).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Scroll the view down and select the cell below the newly added data,
# matching the author's final view state after the edit.
$excel.ActiveWindow.ScrollRow = 43
$ws.Cells.Item($newRow + 1, 3).Select()
